$d = $word.ActiveDocument
$d.Content.Find.Execute("1993-2020", $true, $false, $false, $false, $false, $true, 1, $false, "1993-2024", 2)
